$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (F7) already holds the "FMOD started" status together with its
# "Neutral" cell style (s="8"). Reuse that formatting for F2 and F3 so the
# new status text picks up the matching style instead of creating a
# duplicate style entry.
$ws.Range("F7").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("F3").PasteSpecial(-4122)

# Update the status values: two more sounds have begun work in FMOD.
$ws.Range("F2").Value = "FMOD started"
$ws.Range("F3").Value = "FMOD started"

# Move the active selection to H6, matching the saved view state.
$ws.Range("H6").Select()
